# Auto-generated edit script: updates H:N profit columns per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 990
$ws.Range("I28").Value = 990
$ws.Range("K28").Value = 990
$ws.Range("M28").Value = -505
$ws.Range("H40").Value = 3344.7896
$ws.Range("I40").Value = 2700.15
$ws.Range("J40").Value = 4061.0557
$ws.Range("K40").Value = 2700.15
$ws.Range("L40").Value = 4061.0557
$ws.Range("M40").Value = -2525.15
$ws.Range("N40").Value = -4411.0557
$ws.Range("H55").Value = 1277.0555
$ws.Range("J55").Value = 1807.3334
$ws.Range("L55").Value = 1807.3334
$ws.Range("N55").Value = -2235.3334
$ws.Range("H62").Value = 3956.5
$ws.Range("I62").Value = 3146.6
$ws.Range("K62").Value = 3146.6
$ws.Range("M62").Value = -2522.6
$ws.Range("H65").Value = 3956.5
$ws.Range("I65").Value = 3146.6
$ws.Range("K65").Value = 15733
$ws.Range("M65").Value = -12613
$ws.Range("H107").Value = 2234.7144
$ws.Range("I107").Value = 1598.5
$ws.Range("K107").Value = 1598.5
$ws.Range("M107").Value = 321.5
$ws.Range("H111").Value = 3303
$ws.Range("I111").Value = 2014.5
$ws.Range("K111").Value = 6043.5
$ws.Range("M111").Value = -2976.5
$ws.Range("H138").Value = 3219.1177
$ws.Range("J138").Value = 3421.8667
$ws.Range("L138").Value = 10265.6001
$ws.Range("N138").Value = -20545.6001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 110000
$ws.Range("J117").Value = 110000
$ws.Range("L117").Value = 110000
$ws.Range("N117").Value = -119178
$ws.Range("H122").Value = 1456.579
$ws.Range("I122").Value = 1226.5714
$ws.Range("K122").Value = 3679.7142
$ws.Range("M122").Value = -1229.7142

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1385.1111
$ws.Range("I64").Value = 1371.5
$ws.Range("J64").Value = 1412.3334
$ws.Range("K64").Value = 1371.5
$ws.Range("L64").Value = 1412.3334
$ws.Range("M64").Value = -1146.5
$ws.Range("N64").Value = -1862.3334
$ws.Range("H67").Value = 1385.1111
$ws.Range("I67").Value = 1371.5
$ws.Range("J67").Value = 1412.3334
$ws.Range("K67").Value = 1371.5
$ws.Range("L67").Value = 1412.3334
$ws.Range("M67").Value = -591.5
$ws.Range("N67").Value = -2972.3334
$ws.Range("H110").Value = 114326.336
$ws.Range("J110").Value = 114326.336
$ws.Range("L110").Value = 114326.336
$ws.Range("N110").Value = -122506.336

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = $null
$ws.Range("H64").Value = 108000
$ws.Range("J64").Value = 108000
$ws.Range("L64").Value = 108000
$ws.Range("N64").Value = -108496
$ws.Range("H67").Value = 108000
$ws.Range("J67").Value = 108000
$ws.Range("L67").Value = 108000
$ws.Range("N67").Value = -109716

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2582.9333
$ws.Range("J34").Value = 4374.75
$ws.Range("L34").Value = 13124.25
$ws.Range("N34").Value = -13292.25
$ws.Range("H39").Value = 229401.56
$ws.Range("J39").Value = 219757
$ws.Range("L39").Value = 659271
$ws.Range("N39").Value = -659859
$ws.Range("H55").Value = 11044.444
$ws.Range("I55").Value = 9500
$ws.Range("J55").Value = 11485.714
$ws.Range("K55").Value = 28500
$ws.Range("L55").Value = 34457.142
$ws.Range("M55").Value = -28323
$ws.Range("N55").Value = -34811.142
$ws.Range("H75").Value = 111241590
$ws.Range("J75").Value = 37210730
$ws.Range("L75").Value = 111632190
$ws.Range("N75").Value = -111634186
$ws.Range("H78").Value = 111241590
$ws.Range("J78").Value = 37210730
$ws.Range("L78").Value = 334896570
$ws.Range("N78").Value = -334906554
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = $null
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H97").Value = 2977436
$ws.Range("I97").Value = 7143476.5
$ws.Range("J97").Value = 1692.8572
$ws.Range("K97").Value = 21430429.5
$ws.Range("L97").Value = 5078.571599999999
$ws.Range("M97").Value = -21429933.5
$ws.Range("N97").Value = -6070.571599999999
$ws.Range("H98").Value = 999.4
$ws.Range("I98").Value = 999.5
$ws.Range("K98").Value = 2998.5
$ws.Range("M98").Value = -1500.5
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 4000
$ws.Range("K99").Value = 12000
$ws.Range("M99").Value = -9754
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -49868
$ws.Range("H103").Value = 1452.4
$ws.Range("J103").Value = 2330
$ws.Range("L103").Value = 6990
$ws.Range("N103").Value = -8748
$ws.Range("H104").Value = 4749.8335
$ws.Range("J104").Value = 5000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -20242
$ws.Range("H112").Value = 3304
$ws.Range("I112").Value = 405.66666
$ws.Range("K112").Value = 1216.99998
$ws.Range("M112").Value = -108.9999800000001
$ws.Range("H134").Value = 5076
$ws.Range("I134").Value = 4278.5454
$ws.Range("K134").Value = 12835.6362
$ws.Range("M134").Value = -7765.636200000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5670863.5
$ws.Range("I11").Value = 24030500
$ws.Range("J11").Value = 570964.3
$ws.Range("K11").Value = 24030500
$ws.Range("L11").Value = 570964.3
$ws.Range("M11").Value = -24030361
$ws.Range("N11").Value = -571242.3
$ws.Range("H48").Value = 21000
$ws.Range("J48").Value = 21000
$ws.Range("L48").Value = 21000
$ws.Range("N48").Value = -21970
$ws.Range("H53").Value = 15000
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16262
$ws.Range("H97").Value = 576.5263
$ws.Range("I97").Value = 606.44446
$ws.Range("K97").Value = 606.44446
$ws.Range("M97").Value = -110.44446
$ws.Range("H132").Value = 58831944
$ws.Range("I132").Value = 71437720
$ws.Range("K132").Value = 214313160
$ws.Range("M132").Value = -214310630

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3936.3125
$ws.Range("I40").Value = 1711.2858
$ws.Range("K40").Value = 1711.2858
$ws.Range("M40").Value = -1575.2858
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 30000
$ws.Range("K42").Value = 30000
$ws.Range("M42").Value = -29437
$ws.Range("H43").Value = 2895425.5
$ws.Range("I43").Value = 4039600
$ws.Range("K43").Value = 4039600
$ws.Range("M43").Value = -4039407
$ws.Range("H46").Value = 2539.5454
$ws.Range("I46").Value = 2176.111
$ws.Range("K46").Value = 2176.111
$ws.Range("M46").Value = -1988.111
$ws.Range("H49").Value = 30000
$ws.Range("I49").Value = 30000
$ws.Range("K49").Value = 30000
$ws.Range("M49").Value = -29853
$ws.Range("H68").Value = 4732.6665
$ws.Range("I68").Value = 3574.25
$ws.Range("J68").Value = 14000
$ws.Range("K68").Value = 3574.25
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = -2825.25
$ws.Range("N68").Value = -15498
$ws.Range("H71").Value = 4732.6665
$ws.Range("I71").Value = 3574.25
$ws.Range("J71").Value = 14000
$ws.Range("K71").Value = 17871.25
$ws.Range("L71").Value = 70000
$ws.Range("M71").Value = -14127.25
$ws.Range("N71").Value = -77488

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 41536.668
$ws.Range("I55").Value = 62227.5
$ws.Range("J55").Value = 155
$ws.Range("K55").Value = 62227.5
$ws.Range("L55").Value = 155
$ws.Range("M55").Value = -61950.5
$ws.Range("N55").Value = -709
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21476
$ws.Range("H62").Value = 4451683.5
$ws.Range("J62").Value = 40003990
$ws.Range("L62").Value = 40003990
$ws.Range("N62").Value = -40005238
$ws.Range("H65").Value = 4451683.5
$ws.Range("J65").Value = 40003990
$ws.Range("L65").Value = 200019950
$ws.Range("N65").Value = -200026190
$ws.Range("H75").Value = 9418324
$ws.Range("J75").Value = 9418324
$ws.Range("L75").Value = 9418324
$ws.Range("N75").Value = -9420196
$ws.Range("H78").Value = 9418324
$ws.Range("J78").Value = 9418324
$ws.Range("L78").Value = 28254972
$ws.Range("N78").Value = -28264332
$ws.Range("H107").Value = 22728682
$ws.Range("I107").Value = 31251626
$ws.Range("K107").Value = 93754878
$ws.Range("M107").Value = -93752958
$ws.Range("H122").Value = 6232.8
$ws.Range("I122").Value = 3809.5
$ws.Range("K122").Value = 11428.5
$ws.Range("M122").Value = -8978.5
$ws.Range("H126").Value = 1243.2
$ws.Range("I126").Value = 1197.75
$ws.Range("K126").Value = 3593.25
$ws.Range("M126").Value = -1123.25
